# Auto-generated: refresh market-price / profit figures across the Leve
# profit tracker sheets (currentAveragePrice* / LevePrice* / LeveProfit*),
# as pulled by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 12999.429
$ws.Cells.Item(32, 9).Value = 14999.333
$ws.Cells.Item(32, 11).Value = 14999.333
$ws.Cells.Item(32, 13).Value = -14673.333
$ws.Cells.Item(115, 8).Value = 6995.75
$ws.Cells.Item(115, 9).Value = 7762
$ws.Cells.Item(115, 11).Value = 23286
$ws.Cells.Item(115, 13).Value = -21719
$ws.Cells.Item(137, 8).Value = 2064.2727
$ws.Cells.Item(137, 9).Value = 1821.0588
$ws.Cells.Item(137, 11).Value = 5463.1764
$ws.Cells.Item(137, 13).Value = -2913.1764
$ws.Cells.Item(138, 8).Value = 2636.9412
$ws.Cells.Item(138, 10).Value = 3171.6177
$ws.Cells.Item(138, 12).Value = 9514.8531
$ws.Cells.Item(138, 14).Value = -19794.8531

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1886.6875
$ws.Cells.Item(2, 9).Value = 2256.4167
$ws.Cells.Item(2, 11).Value = 2256.4167
$ws.Cells.Item(2, 13).Value = -2143.4167
$ws.Cells.Item(32, 8).Value = 5293.4546
$ws.Cells.Item(32, 9).Value = 412.58334
$ws.Cells.Item(32, 11).Value = 412.58334
$ws.Cells.Item(32, 13).Value = -125.58334
$ws.Cells.Item(61, 8).Value = 3403.3462
$ws.Cells.Item(61, 9).Value = 3124.5
$ws.Cells.Item(61, 10).Value = 6749.5
$ws.Cells.Item(61, 11).Value = 3124.5
$ws.Cells.Item(61, 12).Value = 6749.5
$ws.Cells.Item(61, 13).Value = -2912.5
$ws.Cells.Item(61, 14).Value = -7173.5
$ws.Cells.Item(112, 8).Value = 20666.334
$ws.Cells.Item(112, 10).Value = 20666.334
$ws.Cells.Item(112, 12).Value = 20666.334
$ws.Cells.Item(112, 14).Value = -23620.334
$ws.Cells.Item(116, 8).Value = 1886.6875
$ws.Cells.Item(116, 9).Value = 2256.4167
$ws.Cells.Item(116, 11).Value = 2256.4167
$ws.Cells.Item(116, 13).Value = 37.58329999999978
$ws.Cells.Item(118, 8).Value = 70000
$ws.Cells.Item(118, 10).Value = 70000
$ws.Cells.Item(118, 12).Value = 70000
$ws.Cells.Item(118, 14).Value = -73314
$ws.Cells.Item(136, 8).Value = 3403.3462
$ws.Cells.Item(136, 9).Value = 3124.5
$ws.Cells.Item(136, 10).Value = 6749.5
$ws.Cells.Item(136, 11).Value = 9373.5
$ws.Cells.Item(136, 12).Value = 20248.5
$ws.Cells.Item(136, 13).Value = -6823.5
$ws.Cells.Item(136, 14).Value = -25348.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1886.6875
$ws.Cells.Item(3, 9).Value = 2256.4167
$ws.Cells.Item(3, 11).Value = 2256.4167
$ws.Cells.Item(3, 13).Value = -2142.4167

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 2008.2
$ws.Cells.Item(7, 9).Value = 2008.2
$ws.Cells.Item(7, 11).Value = 2008.2
$ws.Cells.Item(7, 13).Value = -1895.2
$ws.Cells.Item(16, 8).Value = 19530.562
$ws.Cells.Item(16, 9).Value = 12806.923
$ws.Cells.Item(16, 11).Value = 12806.923
$ws.Cells.Item(16, 13).Value = -12519.923
$ws.Cells.Item(31, 8).Value = 2796.0303
$ws.Cells.Item(31, 9).Value = 2050.4443
$ws.Cells.Item(31, 10).Value = 3075.625
$ws.Cells.Item(31, 11).Value = 2050.4443
$ws.Cells.Item(31, 12).Value = 3075.625
$ws.Cells.Item(31, 13).Value = -1755.4443
$ws.Cells.Item(31, 14).Value = -3665.625
$ws.Cells.Item(34, 8).Value = 2796.0303
$ws.Cells.Item(34, 9).Value = 2050.4443
$ws.Cells.Item(34, 10).Value = 3075.625
$ws.Cells.Item(34, 11).Value = 2050.4443
$ws.Cells.Item(34, 12).Value = 3075.625
$ws.Cells.Item(34, 13).Value = -1848.4443
$ws.Cells.Item(34, 14).Value = -3479.625
$ws.Cells.Item(86, 8).Value = 6988.3076
$ws.Cells.Item(86, 9).Value = 6324.8335
$ws.Cells.Item(86, 10).Value = 7557
$ws.Cells.Item(86, 11).Value = 6324.8335
$ws.Cells.Item(86, 12).Value = 7557
$ws.Cells.Item(86, 13).Value = -5201.8335
$ws.Cells.Item(86, 14).Value = -9803
$ws.Cells.Item(89, 8).Value = 6988.3076
$ws.Cells.Item(89, 9).Value = 6324.8335
$ws.Cells.Item(89, 10).Value = 7557
$ws.Cells.Item(89, 11).Value = 31624.1675
$ws.Cells.Item(89, 12).Value = 37785
$ws.Cells.Item(89, 13).Value = -26008.1675
$ws.Cells.Item(89, 14).Value = -49017
$ws.Cells.Item(99, 8).Value = 5523.75
$ws.Cells.Item(99, 9).Value = 4269.5713
$ws.Cells.Item(99, 11).Value = 4269.5713
$ws.Cells.Item(99, 13).Value = -2771.5713
$ws.Cells.Item(113, 8).Value = 19530.562
$ws.Cells.Item(113, 9).Value = 12806.923
$ws.Cells.Item(113, 11).Value = 12806.923
$ws.Cells.Item(113, 13).Value = -10636.923
$ws.Cells.Item(126, 8).Value = 5523.75
$ws.Cells.Item(126, 9).Value = 4269.5713
$ws.Cells.Item(126, 11).Value = 12808.7139
$ws.Cells.Item(126, 13).Value = -10338.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 374.44446
$ws.Cells.Item(18, 9).Value = 379.625
$ws.Cells.Item(18, 10).Value = 333
$ws.Cells.Item(18, 11).Value = 1138.875
$ws.Cells.Item(18, 12).Value = 999
$ws.Cells.Item(18, 13).Value = -969.875
$ws.Cells.Item(18, 14).Value = -1337
$ws.Cells.Item(122, 8).Value = 6170.8335
$ws.Cells.Item(122, 9).Value = 383.33334
$ws.Cells.Item(122, 10).Value = 11958.333
$ws.Cells.Item(122, 11).Value = 3450.00006
$ws.Cells.Item(122, 12).Value = 107624.997
$ws.Cells.Item(122, 13).Value = -1000.00006
$ws.Cells.Item(122, 14).Value = -112524.997
$ws.Cells.Item(138, 8).Value = 9266.666999999999
$ws.Cells.Item(138, 10).Value = 11650
$ws.Cells.Item(138, 12).Value = 34950
$ws.Cells.Item(138, 14).Value = -45230
$ws.Cells.Item(139, 8).Value = 7525.1177
$ws.Cells.Item(139, 9).Value = 1595
$ws.Cells.Item(139, 11).Value = 4785
$ws.Cells.Item(139, 13).Value = 355
$ws.Cells.Item(140, 8).Value = 2686.875
$ws.Cells.Item(140, 9).Value = 1624.9
$ws.Cells.Item(140, 11).Value = 4874.700000000001
$ws.Cells.Item(140, 13).Value = 305.2999999999993
$ws.Cells.Item(141, 8).Value = 7874
$ws.Cells.Item(141, 9).Value = 3798
$ws.Cells.Item(141, 11).Value = 11394
$ws.Cells.Item(141, 13).Value = -6214

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 261.55554
$ws.Cells.Item(2, 9).Value = 342.16666
$ws.Cells.Item(2, 11).Value = 342.16666
$ws.Cells.Item(2, 13).Value = -229.16666
$ws.Cells.Item(26, 8).Value = 35000
$ws.Cells.Item(26, 10).Value = 35000
$ws.Cells.Item(26, 12).Value = 35000
$ws.Cells.Item(26, 14).Value = -35560
$ws.Cells.Item(47, 8).Value = 25000
$ws.Cells.Item(47, 10).Value = 25000
$ws.Cells.Item(47, 12).Value = 25000
$ws.Cells.Item(47, 14).Value = -26136
$ws.Cells.Item(50, 8).Value = 35000
$ws.Cells.Item(50, 10).Value = 35000
$ws.Cells.Item(50, 12).Value = 35000
$ws.Cells.Item(50, 14).Value = -35996
$ws.Cells.Item(70, 8).Value = 4989.5
$ws.Cells.Item(70, 9).Value = 4989.5
$ws.Cells.Item(70, 11).Value = 4989.5
$ws.Cells.Item(70, 13).Value = -4719.5
$ws.Cells.Item(73, 8).Value = 4989.5
$ws.Cells.Item(73, 9).Value = 4989.5
$ws.Cells.Item(73, 11).Value = 4989.5
$ws.Cells.Item(73, 13).Value = -4053.5
$ws.Cells.Item(102, 8).Value = 2226.8
$ws.Cells.Item(102, 9).Value = 2274.7778
$ws.Cells.Item(102, 11).Value = 2274.7778
$ws.Cells.Item(102, 13).Value = -652.7777999999998
$ws.Cells.Item(126, 8).Value = 6787.636
$ws.Cells.Item(126, 9).Value = 5824
$ws.Cells.Item(126, 10).Value = 8474
$ws.Cells.Item(126, 11).Value = 17472
$ws.Cells.Item(126, 12).Value = 25422
$ws.Cells.Item(126, 13).Value = -15002
$ws.Cells.Item(126, 14).Value = -30362
$ws.Cells.Item(132, 8).Value = 4960.5625
$ws.Cells.Item(132, 9).Value = 4691.433
$ws.Cells.Item(132, 11).Value = 14074.299
$ws.Cells.Item(132, 13).Value = -11544.299

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 19041
$ws.Cells.Item(7, 9).Value = 19041
$ws.Cells.Item(7, 11).Value = 19041
$ws.Cells.Item(7, 13).Value = -18929
$ws.Cells.Item(16, 8).Value = 926.5714
$ws.Cells.Item(16, 9).Value = 926.5714
$ws.Cells.Item(16, 11).Value = 926.5714
$ws.Cells.Item(16, 13).Value = -756.5714
$ws.Cells.Item(61, 8).Value = 3845.6428
$ws.Cells.Item(61, 10).Value = 6335
$ws.Cells.Item(61, 12).Value = 6335
$ws.Cells.Item(61, 14).Value = -6739
$ws.Cells.Item(113, 8).Value = 3845.6428
$ws.Cells.Item(113, 10).Value = 6335
$ws.Cells.Item(113, 12).Value = 6335
$ws.Cells.Item(113, 14).Value = -10675
$ws.Cells.Item(122, 8).Value = 3948.739
$ws.Cells.Item(122, 9).Value = 4134.7
$ws.Cells.Item(122, 10).Value = 2709
$ws.Cells.Item(122, 11).Value = 12404.1
$ws.Cells.Item(122, 12).Value = 8127
$ws.Cells.Item(122, 13).Value = -9954.099999999999
$ws.Cells.Item(122, 14).Value = -13027
$ws.Cells.Item(126, 8).Value = 19041
$ws.Cells.Item(126, 9).Value = 19041
$ws.Cells.Item(126, 11).Value = 57123
$ws.Cells.Item(126, 13).Value = -54653
$ws.Cells.Item(135, 8).Value = 89997.5
$ws.Cells.Item(135, 10).Value = 89997.5
$ws.Cells.Item(135, 12).Value = 89997.5
$ws.Cells.Item(135, 14).Value = -100137.5
$ws.Cells.Item(136, 8).Value = 1185.12
$ws.Cells.Item(136, 9).Value = 996.3182
$ws.Cells.Item(136, 11).Value = 2988.9546
$ws.Cells.Item(136, 13).Value = -438.9546
$ws.Cells.Item(139, 8).Value = 89499
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 89499
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 89499
$ws.Cells.Item(139, 14).Value = -99779
$ws.Cells.Item(139, 13).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3792.8518
$ws.Cells.Item(136, 9).Value = 3615.8462
$ws.Cells.Item(136, 11).Value = 10847.5386
$ws.Cells.Item(136, 13).Value = -8297.5386

